$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# BLEU score (B11)
$ws.Range("B11").Value = 0.09386265332331679

# Code BLEU (B12) + its note (C12)
$ws.Range("B12").Value = 0.2739426453326889
$ws.Range("C12").Value = "{'codebleu': 0.27394264533268886, 'ngram_match_score': 0.09352561860594595, 'weighted_ngram_match_score': 0.17527767713683412, 'syntax_match_score': 0.45517241379310347, 'dataflow_match_score': 0.3717948717948718}"

# Embeddings and Cosine similarity (B13)
$ws.Range("B13").Value = 0.5961298176503818
